# Updates the crypto price/volume table to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.808.74'
$ws.Cells.Item(2, 5).Value = '  -0.32%  '

$ws.Cells.Item(3, 4).Value = '1.888.70'
$ws.Cells.Item(3, 5).Value = '  -0.80%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  +0.11%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.7909'
$ws.Cells.Item(5, 5).Value = '  -1.06%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '242.01'
$ws.Cells.Item(6, 5).Value = '  +0.56%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.002'
$ws.Cells.Item(7, 5).Value = '  +0.14%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3174'
$ws.Cells.Item(8, 5).Value = '  +2.20%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '25.42'
$ws.Cells.Item(9, 5).Value = '  -3.56%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.07042'
$ws.Cells.Item(10, 5).Value = '  +1.12%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08054'
$ws.Cells.Item(11, 5).Value = '  +0.93%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.7676'
$ws.Cells.Item(12, 5).Value = '  +4.11%  '

$ws.Cells.Item(13, 4).Value = '1.869.80'
$ws.Cells.Item(13, 5).Value = '  -2.86%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.296'
$ws.Cells.Item(14, 5).Value = '  +2.63%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '92.06'
$ws.Cells.Item(15, 5).Value = '  -0.17%  '

$ws.Cells.Item(16, 4).Value = '29.836.10'
$ws.Cells.Item(16, 5).Value = '  -0.28%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '13.82'
$ws.Cells.Item(17, 5).Value = '  -0.76%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '5.921'
$ws.Cells.Item(18, 5).Value = '  +1.33%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '243.48'
$ws.Cells.Item(19, 5).Value = '  -0.28%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000007716'
$ws.Cells.Item(20, 5).Value = '  -0.34%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.002'
$ws.Cells.Item(21, 5).Value = '  +0.10%  '

$ws.Cells.Item(22, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(22, 4).Value = '2.140.59'
$ws.Cells.Item(22, 5).Value = '  -0.55%  '

$ws.Cells.Item(23, 2).Value = 'Chainlink'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '8.066'
$ws.Cells.Item(23, 5).Value = '  +16.70%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.002'
$ws.Cells.Item(24, 5).Value = '  +0.13%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.1619'
$ws.Cells.Item(25, 5).Value = '  +11.97%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.299'
$ws.Cells.Item(26, 5).Value = '  +1.27%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '163.85'
$ws.Cells.Item(27, 5).Value = '  -2.11%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.67'
$ws.Cells.Item(28, 5).Value = '  -0.86%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.054'
$ws.Cells.Item(29, 5).Value = '  -0.03%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.370'
$ws.Cells.Item(30, 5).Value = '  +1.17%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.537'
$ws.Cells.Item(31, 5).Value = '  +1.77%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.454'
$ws.Cells.Item(32, 5).Value = '  +4.17%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05653'
$ws.Cells.Item(33, 5).Value = '  +2.38%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.104'
$ws.Cells.Item(34, 5).Value = '  +1.30%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.262'
$ws.Cells.Item(35, 5).Value = '  -0.05%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7346'
$ws.Cells.Item(36, 5).Value = '  +0.68%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.9984'
$ws.Cells.Item(37, 5).Value = '  -0.10%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.715'
$ws.Cells.Item(38, 5).Value = '  -0.19%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01926'
$ws.Cells.Item(39, 5).Value = '  +0.25%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.777'
$ws.Cells.Item(40, 5).Value = '  -0.29%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.4420'
$ws.Cells.Item(41, 5).Value = '  +0.66%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '72.37'
$ws.Cells.Item(42, 5).Value = '  +0.54%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.838'
$ws.Cells.Item(43, 5).Value = '  -2.15%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.8416'
$ws.Cells.Item(44, 5).Value = '  +0.70%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.001'
$ws.Cells.Item(45, 5).Value = '  +0.16%  '

$ws.Cells.Item(46, 4).Value = '1.025.34'
$ws.Cells.Item(46, 5).Value = '  +4.78%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.872'
$ws.Cells.Item(47, 5).Value = '  -0.36%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '101.94'
$ws.Cells.Item(48, 5).Value = '  +1.20%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.945'
$ws.Cells.Item(49, 5).Value = '  +2.82%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.456'
$ws.Cells.Item(50, 5).Value = '  -0.99%  '

$ws.Cells.Item(51, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(51, 4).Value = '2.054.78'
$ws.Cells.Item(51, 5).Value = '  -0.61%  '
